# Commit: [minor] - added distinct() in q3_rdd
# The underlying data change is that the RDD value for Q2 (cell B4) is
# updated from 75.44 to 73.75. The active selection is also moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RDD / Q2 value (currently 75.44 -> 73.75). The figures in
# this table are stored as text (shared strings), not numbers, so a
# direct Value assignment of "73.75" must be prevented from being
# auto-converted to a numeric cell by Excel. Build the text via a
# scratch-cell formula, then paste only the resulting value back into
# B4; this preserves the text cell type without touching the cell's
# number format / style (which a NumberFormat="@" assignment would do).
$ws.Range("ZZ1").Formula = '=TEXT(73.75,"0.00")'
$ws.Range("ZZ1").Copy()
$ws.Range("B4").PasteSpecial(-4163) # xlPasteValues
$ws.Range("ZZ1").EntireColumn.Delete()
$excel.CutCopyMode = 0

# Move the active selection to J13, matching the recorded sheet view state
$ws.Range("J13").Select()

$wb.Save()
